{"js": "// Update the date heading and every \"A\u00d7B=C\" multiplication-table answer\n// cell from the previous day's values to the new day's values. Every\n// <w:t> run in the document changes exactly once, and every old value is\n// unique, so we can safely locate each run via body.search() and replace\n// its text in place (preserving the run's formatting).\nconst replacements = [\n  [\"2025-06-30 Monday\", \"2025-07-01 Tuesday\"],\n  [\"508\u00d74=2032\", \"553\u00d76=3318\"],\n  [\"398\u00d79=3582\", \"634\u00d79=5706\"],\n  [\"967\u00d74=3868\", \"845\u00d74=3380\"],\n  [\"225\u00d72=450\", \"395\u00d73=1185\"],\n  [\"219\u00d72=438\", \"119\u00d78=952\"],\n  [\"545\u00d73=1635\", \"736\u00d77=5152\"],\n  [\"794\u00d74=3176\", \"220\u00d79=1980\"],\n  [\"317\u00d76=1902\", \"508\u00d77=3556\"],\n  [\"360\u00d75=1800\", \"582\u00d73=1746\"],\n  [\"993\u00d79=8937\", \"699\u00d72=1398\"],\n  [\"442\u00d79=3978\", \"662\u00d76=3972\"],\n  [\"257\u00d75=1285\", \"960\u00d74=3840\"],\n  [\"133\u00d72=266\", \"436\u00d75=2180\"],\n  [\"322\u00d75=1610\", \"876\u00d77=6132\"],\n  [\"728\u00d76=4368\", \"508\u00d76=3048\"],\n  [\"846\u00d72=1692\", \"308\u00d78=2464\"],\n  [\"162\u00d73=486\", \"529\u00d79=4761\"],\n  [\"682\u00d73=2046\", \"624\u00d75=3120\"],\n  [\"920\u00d76=5520\", \"303\u00d73=909\"],\n  [\"992\u00d74=3968\", \"169\u00d74=676\"],\n  [\"621\u00d76=3726\", \"348\u00d77=2436\"],\n  [\"589\u00d76=3534\", \"799\u00d75=3995\"],\n  [\"134\u00d72=268\", \"246\u00d72=492\"],\n  [\"962\u00d75=4810\", \"859\u00d73=2577\"],\n  [\"321\u00d77=2247\", \"793\u00d72=1586\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date heading and every \"A\u00d7B=C\" multiplication-table answer\n# cell from the previous day's values to the new day's values. Every\n# old value in the document is unique, so Find/Replace against the full\n# document range can locate and update each run safely, one at a time.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-30 Monday\", \"2025-07-01 Tuesday\"),\n    @(\"508\u00d74=2032\", \"553\u00d76=3318\"),\n    @(\"398\u00d79=3582\", \"634\u00d79=5706\"),\n    @(\"967\u00d74=3868\", \"845\u00d74=3380\"),\n    @(\"225\u00d72=450\", \"395\u00d73=1185\"),\n    @(\"219\u00d72=438\", \"119\u00d78=952\"),\n    @(\"545\u00d73=1635\", \"736\u00d77=5152\"),\n    @(\"794\u00d74=3176\", \"220\u00d79=1980\"),\n    @(\"317\u00d76=1902\", \"508\u00d77=3556\"),\n    @(\"360\u00d75=1800\", \"582\u00d73=1746\"),\n    @(\"993\u00d79=8937\", \"699\u00d72=1398\"),\n    @(\"442\u00d79=3978\", \"662\u00d76=3972\"),\n    @(\"257\u00d75=1285\", \"960\u00d74=3840\"),\n    @(\"133\u00d72=266\", \"436\u00d75=2180\"),\n    @(\"322\u00d75=1610\", \"876\u00d77=6132\"),\n    @(\"728\u00d76=4368\", \"508\u00d76=3048\"),\n    @(\"846\u00d72=1692\", \"308\u00d78=2464\"),\n    @(\"162\u00d73=486\", \"529\u00d79=4761\"),\n    @(\"682\u00d73=2046\", \"624\u00d75=3120\"),\n    @(\"920\u00d76=5520\", \"303\u00d73=909\"),\n    @(\"992\u00d74=3968\", \"169\u00d74=676\"),\n    @(\"621\u00d76=3726\", \"348\u00d77=2436\"),\n    @(\"589\u00d76=3534\", \"799\u00d75=3995\"),\n    @(\"134\u00d72=268\", \"246\u00d72=492\"),\n    @(\"962\u00d75=4810\", \"859\u00d73=2577\"),\n    @(\"321\u00d77=2247\", \"793\u00d72=1586\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
